# Merge the separate "An" / " " / "image" runs in the caption textbox
# into a single run reading "An image", per the golden-test diff.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Shape 3 is "TextBox 3", the caption under the picture.
$shape = $s.Shapes.Item(3)
$tr = $shape.TextFrame.TextRange

# Re-assigning Text on the whole TextRange is a no-op when the
# concatenation of the existing runs already equals the new string (the
# engine only rewrites runs that actually changed). Go through
# Characters(), covering the full span, so the run-merge is forced even
# though the visible text itself doesn't change.
$c = $tr.Characters(1, $tr.Length)
$c.Text = "An image"
